$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 19.12497595493468
$ws.Range("F2").Value = 1.830708315459377
$ws.Range("E3").Value = 8.0585519296698
$ws.Range("F3").Value = 2.855492514199263
$ws.Range("E4").Value = 8.0585519296698
$ws.Range("F4").Value = 2.855492514199263
$ws.Range("E5").Value = 8.0585519296698
$ws.Range("F5").Value = 2.855492514199263
$ws.Range("E6").Value = 8.0585519296698
$ws.Range("F6").Value = 2.855492514199263
$ws.Range("E7").Value = 10.34407011703182
$ws.Range("F7").Value = 1.016780692257745
$ws.Range("E8").Value = 10.34407011703182
$ws.Range("F8").Value = 1.016780692257745
$ws.Range("E9").Value = 10.34407011703182
$ws.Range("F9").Value = 1.016780692257745
$ws.Range("E10").Value = 10.34407011703182
$ws.Range("F10").Value = 1.016780692257745
$ws.Range("E11").Value = 9.201311023350812
$ws.Range("F11").Value = 1.936136603228504
$ws.Range("E12").Value = 9.201311023350812
$ws.Range("F12").Value = 1.936136603228504
$ws.Range("E13").Value = 9.201311023350812
$ws.Range("F13").Value = 1.936136603228504
$ws.Range("E14").Value = 9.201311023350812
$ws.Range("F14").Value = 1.936136603228504
$ws.Range("E15").Value = 9.201311023350812
$ws.Range("F15").Value = 1.936136603228504
$ws.Range("E16").Value = 9.201311023350812
$ws.Range("F16").Value = 1.936136603228504
$ws.Range("E17").Value = 12.62958830439384
$ws.Range("F17").Value = 2.855492514199263
$ws.Range("E18").Value = 12.62958830439384
$ws.Range("F18").Value = 2.855492514199263
$ws.Range("E19").Value = 12.62958830439384
$ws.Range("F19").Value = 2.855492514199263
$ws.Range("E20").Value = 12.62958830439384
$ws.Range("F20").Value = 2.855492514199263
$ws.Range("E21").Value = 14.91510649175587
$ws.Range("F21").Value = 1.016780692257745
$ws.Range("E22").Value = 14.91510649175587
$ws.Range("F22").Value = 1.016780692257745
$ws.Range("E23").Value = 14.91510649175587
$ws.Range("F23").Value = 1.016780692257745
$ws.Range("E24").Value = 14.91510649175587
$ws.Range("F24").Value = 1.016780692257745
$ws.Range("E25").Value = 13.77234739807485
$ws.Range("F25").Value = 1.936136603228504
$ws.Range("E26").Value = 13.77234739807485
$ws.Range("F26").Value = 1.936136603228504
$ws.Range("E27").Value = 13.77234739807485
$ws.Range("F27").Value = 1.936136603228504
$ws.Range("E28").Value = 13.77234739807485
$ws.Range("F28").Value = 1.936136603228504
$ws.Range("E29").Value = 13.77234739807485
$ws.Range("F29").Value = 1.936136603228504
$ws.Range("E30").Value = 13.77234739807485
$ws.Range("F30").Value = 1.936136603228504
$ws.Range("E31").Value = 1.644298635518413
$ws.Range("F31").Value = 2.841876901740046
$ws.Range("E32").Value = 1.644298635518413
$ws.Range("F32").Value = 2.841876901740046
$ws.Range("E33").Value = 1.644298635518413
$ws.Range("F33").Value = 2.841876901740046
$ws.Range("E34").Value = 1.644298635518413
$ws.Range("F34").Value = 2.841876901740046
$ws.Range("E35").Value = 3.807437837578521
$ws.Range("F35").Value = 1.054636323345442
$ws.Range("E36").Value = 3.807437837578521
$ws.Range("F36").Value = 1.054636323345442
$ws.Range("E37").Value = 3.807437837578521
$ws.Range("F37").Value = 1.054636323345442
$ws.Range("E38").Value = 3.807437837578521
$ws.Range("F38").Value = 1.054636323345442
$ws.Range("E39").Value = 2.725868236548467
$ws.Range("F39").Value = 1.948256612542744
$ws.Range("E40").Value = 2.725868236548467
$ws.Range("F40").Value = 1.948256612542744
$ws.Range("E41").Value = 2.725868236548467
$ws.Range("F41").Value = 1.948256612542744
$ws.Range("E42").Value = 2.725868236548467
$ws.Range("F42").Value = 1.948256612542744
$ws.Range("E43").Value = 2.725868236548467
$ws.Range("F43").Value = 1.948256612542744
$ws.Range("E44").Value = 2.725868236548467
$ws.Range("F44").Value = 1.948256612542744
$ws.Range("E45").Value = 5.811207670050094
$ws.Range("F45").Value = 1.858931006443455
$ws.Range("E46").Value = 5.811207670050094
$ws.Range("F46").Value = 1.858931006443455
$ws.Range("E47").Value = 0.2863832658242556
$ws.Range("F47").Value = 1.892266529044012
$ws.Range("E48").Value = 0.2863832658242556
$ws.Range("F48").Value = 1.892266529044012
$ws.Range("E49").Value = 0.2863832658242556
$ws.Range("F49").Value = 1.892266529044012
$ws.Range("E50").Value = 0.2863832658242556
$ws.Range("F50").Value = 1.892266529044012